$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Write the six brand-new names first, in the order they were originally
# entered (this also controls the order new entries land in the shared
# string table), at the final alphabetically-sorted row each belongs to.
$ws.Range("A11").Value = "Leandro Freire da Silva"
$ws.Range("A8").Value  = "Jonatas Pereira Cabral de Araujo"
$ws.Range("A31").Value = "Ronaldo da Costa Tavares"
$ws.Range("A25").Value = "Patrícia Pedrosa Alves Braga "
$ws.Range("A13").Value = "Lídice F do Carmo dos Santos"
$ws.Range("A9").Value  = "Julia Valerio Andrade"

# Fill in the rest of the alphabetised roster (A3:A38), reusing the
# existing names/shared strings.
$ws.Range("A3").Value  = "Alane Dantas de Azevedo Lima"
$ws.Range("A4").Value  = "Bruno Roberto Santana Bello"
$ws.Range("A5").Value  = "Fabio Mortari"
$ws.Range("A6").Value  = "Fabrício Kassardjian"
$ws.Range("A7").Value  = "Igor Camargo Garcia"
$ws.Range("A10").Value = "Juliana de Carvalho Fernandes"
$ws.Range("A12").Value = "Leda Miranda Guimarães"
$ws.Range("A14").Value = "Lincoln Stuart Lima da Gama"
$ws.Range("A15").Value = "Luan Morais de Brito"
$ws.Range("A16").Value = "Luciana de Barros Valentino"
$ws.Range("A17").Value = "Luis Felipe Maior"
$ws.Range("A18").Value = "Marcio Trindade"
$ws.Range("A19").Value = "Marcos Renann Fernandes da Silva"
$ws.Range("A20").Value = "Mariana Oliveira Campos Machado"
$ws.Range("A21").Value = "Mateus Menezes Ribeiro"
$ws.Range("A22").Value = "Monick Hellen Nogueira Macena"
$ws.Range("A23").Value = "Nadia Ligia Costa dos Santos"
$ws.Range("A24").Value = "Nivaldo Mariano de Carvalho Junior"
$ws.Range("A26").Value = "Pedro Freitas"
$ws.Range("A27").Value = "Pedro Henrique Alves Rosendo"
$ws.Range("A28").Value = "Pedro Ivan Chaves Oliveira"
$ws.Range("A29").Value = "Pollyana Gomes Minatel"
$ws.Range("A30").Value = "Roberto Freixeira da Silva Junior"
$ws.Range("A32").Value = "Thais  Brasil Lenhard"
$ws.Range("A33").Value = "Thiago de Oliveira Dutra"
$ws.Range("A34").Value = "Vanessa Maria Ramos Fischer"
$ws.Range("A35").Value = "Vanessa Sharine Careaga Camelo"
$ws.Range("A36").Value = "Vinicios Alves de Andrade"
$ws.Range("A37").Value = "Vinicius Padovan Trapia"
$ws.Range("A38").Value = "Vitor Soares Silva"

# Clear the now-stale "Grupo 1" marks that used to sit on rows 25/26 before
# the roster was re-sorted.
$ws.Range("B25").Clear()
$ws.Range("B26").Clear()

# "Grupo 1" ("X") markers stay attached to the same three students (Bruno,
# Pollyana and Roberto), who are now on rows 4, 29 and 30 after the re-sort.
$ws.Range("B4").Value = "X"
$ws.Range("B29").Value = "X"
$ws.Range("B30").Value = "X"

# Reflect the saved selection/view state.
$ws.Range("A13").Select()
